# fix: excel files updated
# Update the WhatsApp Number value in E2 and move the active selection to H5,
# matching the author's committed change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the phone number value stored in E2 (numeric cell, not a string)
$ws.Range("E2").Value = 8731903147

# Move the active cell / selection to H5 as recorded in the sheetView
$ws.Range("H5").Select() | Out-Null
